$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: main heading and bold summary line at the end)
Replace-Text "Play Double Lucky Mushrooms DoubleMax for Free - Slot Game Review" "Play Double Lucky Mushrooms DoubleMax Free - Review"

# "What we like" bullet list
Replace-Text "Beautiful graphics that transport the player to the endless green fields of Ireland" "Beautiful graphics with a natural and relaxing atmosphere"
Replace-Text "Large symbols make it easy to understand when a combo has been won" "Easy to understand when a combo has been won"
Replace-Text "Automatic spins feature with loss limits" "Automatic spins with the option to set limits on losses"
Replace-Text "High volatility with a theoretical RTP of 95.5%" "Demo version available to try out the game before betting with real money"

# "What we don't like" bullet list
Replace-Text "Demo version may not fully represent the game's true potential" "High volatility may not appeal to all players"

# Closing italic summary paragraph
Replace-Text "Experience the magic of Double Lucky Mushrooms DoubleMax with beautiful graphics, high volatility, and a theoretical RTP of 95.5%. Try it out for free now!" "Experience the magic of Double Lucky Mushrooms DoubleMax with beautiful graphics and high volatility."
